# Convert column A (Timestamp) from inline-string text to true Excel date/time
# serial values, formatted with a custom "YYYY-MM-DD HH:MM:SS" number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> Excel date serial (days since 1899-12-30) for each of the 50 data rows.
$serials = [ordered]@{
  2 = 45978.37517361111  # 2025-11-17 09:00:15
  3 = 45978.375231481485  # 2025-11-17 09:00:20
  4 = 45978.37540509259  # 2025-11-17 09:00:35
  5 = 45978.3755787037  # 2025-11-17 09:00:50
  6 = 45978.37569444445  # 2025-11-17 09:01:00
  7 = 45978.37615740741  # 2025-11-17 09:01:40
  8 = 45978.37621527778  # 2025-11-17 09:01:45
  9 = 45978.37633101852  # 2025-11-17 09:01:55
  10 = 45978.37650462963  # 2025-11-17 09:02:10
  11 = 45978.37662037037  # 2025-11-17 09:02:20
  12 = 45978.376967592594  # 2025-11-17 09:02:50
  13 = 45978.37719907407  # 2025-11-17 09:03:10
  14 = 45978.37737268519  # 2025-11-17 09:03:25
  15 = 45978.377546296295  # 2025-11-17 09:03:40
  16 = 45978.37771990741  # 2025-11-17 09:03:55
  17 = 45978.37789351852  # 2025-11-17 09:04:10
  18 = 45978.378067129626  # 2025-11-17 09:04:25
  19 = 45978.378125  # 2025-11-17 09:04:30
  20 = 45978.37829861111  # 2025-11-17 09:04:45
  21 = 45978.37841435185  # 2025-11-17 09:04:55
  22 = 45978.378530092596  # 2025-11-17 09:05:05
  23 = 45978.37876157407  # 2025-11-17 09:05:25
  24 = 45978.37899305556  # 2025-11-17 09:05:45
  25 = 45978.37945601852  # 2025-11-17 09:06:25
  26 = 45978.37957175926  # 2025-11-17 09:06:35
  27 = 45978.37962962963  # 2025-11-17 09:06:40
  28 = 45978.38009259259  # 2025-11-17 09:07:20
  29 = 45978.380208333336  # 2025-11-17 09:07:30
  30 = 45978.380324074074  # 2025-11-17 09:07:40
  31 = 45978.38055555556  # 2025-11-17 09:08:00
  32 = 45978.38078703704  # 2025-11-17 09:08:20
  33 = 45978.38101851852  # 2025-11-17 09:08:40
  34 = 45978.38148148148  # 2025-11-17 09:09:20
  35 = 45978.38159722222  # 2025-11-17 09:09:30
  36 = 45978.38171296296  # 2025-11-17 09:09:40
  37 = 45978.381886574076  # 2025-11-17 09:09:55
  38 = 45978.38234953704  # 2025-11-17 09:10:35
  39 = 45978.38269675926  # 2025-11-17 09:11:05
  40 = 45978.3828125  # 2025-11-17 09:11:15
  41 = 45978.38287037037  # 2025-11-17 09:11:20
  42 = 45978.38321759259  # 2025-11-17 09:11:50
  43 = 45978.38333333333  # 2025-11-17 09:12:00
  44 = 45978.38379629629  # 2025-11-17 09:12:40
  45 = 45978.38391203704  # 2025-11-17 09:12:50
  46 = 45978.38396990741  # 2025-11-17 09:12:55
  47 = 45978.38408564815  # 2025-11-17 09:13:05
  48 = 45978.38425925926  # 2025-11-17 09:13:20
  49 = 45978.38460648148  # 2025-11-17 09:13:50
  50 = 45978.384722222225  # 2025-11-17 09:14:00
  51 = 45978.384780092594  # 2025-11-17 09:14:05
}

$first = $true
foreach ($row in $serials.Keys) {
  $cell = $ws.Cells.Item($row, 1)
  $cell.Value = $serials[$row]
  if ($first) {
    # Establish both the lowercase and the final uppercase custom date/time
    # format codes in the style table (numFmtId 164 then 165), matching how the
    # format was iterated on for the very first cell touched.
    $cell.NumberFormat = "yyyy-mm-dd h:mm:ss"
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $first = $false
  } else {
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
  }
}
